$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 18516
$ws.Range("E2").Value = -357
$ws.Range("F2").Value = -357
$ws.Range("G2").Value = -903
$ws.Range("H2").Value = -864
$ws.Range("I2").Value = -583
$ws.Range("J2").Value = -282
$ws.Range("K2").Value = 10938
$ws.Range("L2").Value = 7200
$ws.Range("M2").Value = 3738
$ws.Range("N2").Value = 3341
$ws.Range("O2").Value = 397
$ws.Range("P2").Value = 764
$ws.Range("Q2").Value = -79
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = -845
$ws.Range("T2").Value = 120
$ws.Range("U2").Value = -200
$ws.Range("V2").Value = 5440
$ws.Range("W2").Value = -1.93
$ws.Range("X2").Value = -4.67
$ws.Range("Y2").Value = -16.06
$ws.Range("Z2").Value = -7.25
$ws.Range("AA2").Value = 192.64
$ws.Range("AB2").Value = 355.91
$ws.Range("AC2").Value = -3813
$ws.Range("AD2").Value = -2.39
$ws.Range("AE2").Value = 22598
$ws.Range("AF2").Value = 0.4
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.1
$ws.Range("AI2").Value = -2.54
$ws.Range("AJ2").Value = 15278000

# Row 3
$ws.Range("D3").Value = 14712
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = -426
$ws.Range("H3").Value = -546
$ws.Range("I3").Value = -434
$ws.Range("J3").Value = -112
$ws.Range("K3").Value = 10286
$ws.Range("L3").Value = 6777
$ws.Range("M3").Value = 3508
$ws.Range("N3").Value = 2961
$ws.Range("O3").Value = 547
$ws.Range("P3").Value = 764
$ws.Range("Q3").Value = 565
$ws.Range("R3").Value = -170
$ws.Range("S3").Value = -647
$ws.Range("T3").Value = 139
$ws.Range("U3").Value = 426
$ws.Range("V3").Value = 4510
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = -3.71
$ws.Range("Y3").Value = -13.77
$ws.Range("Z3").Value = -5.15
$ws.Range("AA3").Value = 193.19
$ws.Range("AB3").Value = 298.88
$ws.Range("AC3").Value = -2840
$ws.Range("AD3").Value = -3.63
$ws.Range("AE3").Value = 20285
$ws.Range("AF3").Value = 0.51
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 2.91
$ws.Range("AI3").Value = -10.09
$ws.Range("AJ3").Value = 15278000

# Row 4
$ws.Range("D4").Value = 15154
$ws.Range("E4").Value = 626
$ws.Range("F4").Value = 626
$ws.Range("G4").Value = 182
$ws.Range("H4").Value = 79
$ws.Range("I4").Value = 107
$ws.Range("J4").Value = -28
$ws.Range("K4").Value = 10979
$ws.Range("L4").Value = 7538
$ws.Range("M4").Value = 3441
$ws.Range("N4").Value = 2938
$ws.Range("O4").Value = 503
$ws.Range("P4").Value = 764
$ws.Range("Q4").Value = 824
$ws.Range("R4").Value = -97
$ws.Range("S4").Value = -252
$ws.Range("T4").Value = 144
$ws.Range("U4").Value = 680
$ws.Range("V4").Value = 4386
$ws.Range("W4").Value = 4.13
$ws.Range("X4").Value = 0.52
$ws.Range("Y4").Value = 3.62
$ws.Range("Z4").Value = 0.74
$ws.Range("AA4").Value = 219.04
$ws.Range("AB4").Value = 300.73
$ws.Range("AC4").Value = 699
$ws.Range("AD4").Value = 23.31
$ws.Range("AE4").Value = 20230
$ws.Range("AF4").Value = 0.8100000000000001
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 3.07
$ws.Range("AI4").Value = 67.98
$ws.Range("AJ4").Value = 15278000

# Row 5
$ws.Range("D5").Value = 15078
$ws.Range("E5").Value = 311
$ws.Range("F5").Value = 311
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = -9
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = -44
$ws.Range("K5").Value = 9995
$ws.Range("L5").Value = 6554
$ws.Range("M5").Value = 3442
$ws.Range("N5").Value = 2816
$ws.Range("O5").Value = 626
$ws.Range("P5").Value = 764
$ws.Range("Q5").Value = 102
$ws.Range("R5").Value = -462
$ws.Range("S5").Value = 419
$ws.Range("T5").Value = 224
$ws.Range("U5").Value = -122
$ws.Range("V5").Value = 4692
$ws.Range("W5").Value = 2.06
$ws.Range("X5").Value = -0.06
$ws.Range("Y5").Value = 1.21
$ws.Range("Z5").Value = -0.09
$ws.Range("AA5").Value = 190.44
$ws.Range("AB5").Value = 294.65
$ws.Range("AC5").Value = 227
$ws.Range("AD5").Value = 67.06
$ws.Range("AE5").Value = 19387
$ws.Range("AF5").Value = 0.79
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 3.28
$ws.Range("AI5").Value = 209.03
$ws.Range("AJ5").Value = 15278000

# Row 6
$ws.Range("D6").Value = 15313
$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 65
$ws.Range("G6").Value = 98
$ws.Range("H6").Value = -16
$ws.Range("I6").Value = 68
$ws.Range("K6").Value = 9357
$ws.Range("L6").Value = 5332
$ws.Range("M6").Value = 4024
$ws.Range("N6").Value = 3086
$ws.Range("P6").Value = 764
$ws.Range("Q6").Value = -391
$ws.Range("R6").Value = 853
$ws.Range("S6").Value = -867
$ws.Range("T6").Value = 266
$ws.Range("U6").Value = -657
$ws.Range("V6").Value = 3440
$ws.Range("W6").Value = 0.42
$ws.Range("X6").Value = -0.11
$ws.Range("Y6").Value = 2.31
$ws.Range("Z6").Value = -0.17
$ws.Range("AA6").Value = 132.51
$ws.Range("AB6").Value = 247.15
$ws.Range("AC6").Value = 446
$ws.Range("AD6").Value = 22.43
$ws.Range("AE6").Value = 21247
$ws.Range("AF6").Value = 0.47
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 4
$ws.Range("AI6").Value = 85.31
$ws.Range("AJ6").Value = 15278000

# Row 7
$ws.Range("D7").Value = 15309
$ws.Range("E7").Value = 258
$ws.Range("G7").Value = 152
$ws.Range("H7").Value = 128
$ws.Range("I7").Value = 199
$ws.Range("K7").Value = 10524
$ws.Range("L7").Value = 6532
$ws.Range("M7").Value = 3992
$ws.Range("P7").Value = 764
$ws.Range("Q7").Value = -366
$ws.Range("R7").Value = -538
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 178
$ws.Range("U7").Value = -1808
$ws.Range("W7").Value = 1.69
$ws.Range("X7").Value = 0.84
$ws.Range("Z7").Value = 1.29
$ws.Range("AA7").Value = 163.63
$ws.Range("AC7").Value = 1303
$ws.Range("AD7").Value = 7.05
$ws.Range("AG7").Value = 500
$ws.Range("AH7").Value = 5.45
$ws.Range("AI7").Value = 38.39
$ws.Range("N7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()

# Row 8
$ws.Range("D8").Value = 18435
$ws.Range("E8").Value = 427
$ws.Range("G8").Value = 396
$ws.Range("H8").Value = 297
$ws.Range("I8").Value = 337
$ws.Range("K8").Value = 11343
$ws.Range("L8").Value = 7091
$ws.Range("M8").Value = 4252
$ws.Range("P8").Value = 764
$ws.Range("Q8").Value = 301
$ws.Range("R8").Value = -346
$ws.Range("S8").Value = 424
$ws.Range("T8").Value = 186
$ws.Range("U8").Value = 83
$ws.Range("W8").Value = 2.32
$ws.Range("X8").Value = 1.61
$ws.Range("Z8").Value = 2.72
$ws.Range("AA8").Value = 166.77
$ws.Range("AC8").Value = 2206
$ws.Range("AD8").Value = 4.16
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 5.45
$ws.Range("AI8").Value = 22.67
$ws.Range("N8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()

# Row 9
$ws.Range("D9").Value = 18554
$ws.Range("E9").Value = 480
$ws.Range("G9").Value = 450
$ws.Range("H9").Value = 337
$ws.Range("I9").Value = 377
$ws.Range("K9").Value = 11775
$ws.Range("L9").Value = 7222
$ws.Range("M9").Value = 4553
$ws.Range("P9").Value = 764
$ws.Range("Q9").Value = 552
$ws.Range("R9").Value = -352
$ws.Range("S9").Value = 24
$ws.Range("T9").Value = 186
$ws.Range("U9").Value = 334
$ws.Range("W9").Value = 2.59
$ws.Range("X9").Value = 1.82
$ws.Range("Z9").Value = 2.92
$ws.Range("AA9").Value = 158.62
$ws.Range("AC9").Value = 2468
$ws.Range("AD9").Value = 3.72
$ws.Range("AG9").Value = 500
$ws.Range("AH9").Value = 5.45
$ws.Range("AI9").Value = 20.26
$ws.Range("N9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
